function Add-AppleSystemRun {
    param($doc, $pos, [string]$text)

    # Insert the raw text at a zero-length range sitting at $pos (i.e. right
    # before the paragraph mark / after whatever run currently ends there).
    $ins = $doc.Range($pos, $pos)
    $ins.InsertAfter($text)

    $newEnd = $pos + $text.Length
    $newRange = $doc.Range($pos, $newEnd)

    # Re-run Find/Replace scoped to just the text we inserted so it gets
    # stamped with its own rPr (AppleSystemUIFont) without touching the
    # run(s) that precede it.
    $f = $newRange.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $f.Replacement.Font.Name = "AppleSystemUIFont"
    $f.Replacement.Font.NameBi = "AppleSystemUIFont"
    [void]$f.Execute($text, $true, $false, $false, $false, $false, $true, 1, $true, $text, 2)

    return $newEnd
}

$d = $word.ActiveDocument

# --- "Marcus:" answer under Spørsmål 1 --------------------------------
$p = $d.Paragraphs.Item(16)
$r = $p.Range
$pos = $r.End - 1
$pos = Add-AppleSystemRun $d $pos " "
$pos = Add-AppleSystemRun $d $pos "Jeg mener det mest relevante og nyttige er css og html, begge er sentrale og interessante å lære om."

# --- "Marcus:" answer under Spørsmål 2 --------------------------------
$p = $d.Paragraphs.Item(38)
$r = $p.Range
$pos = $r.End - 1
$pos = Add-AppleSystemRun $d $pos " "
$pos = Add-AppleSystemRun $d $pos "Jeg jobber fortsatt med å finne den beste læringsformen for meg, men syns det har vært en god løsning og gått igjennom ting så latt oss jobbe litt med det så gå gjennom det etterpå."

# --- "Marcus:" answer under Spørsmål 3 --------------------------------
$p = $d.Paragraphs.Item(65)
$r = $p.Range
$pos = $r.End - 1
$pos = Add-AppleSystemRun $d $pos " "
$pos = Add-AppleSystemRun $d $pos "Kanskje animasjoner"
$pos = Add-AppleSystemRun $d $pos "."
$pos = Add-AppleSystemRun $d $pos " "
$pos = Add-AppleSystemRun $d $pos "J"
$pos = Add-AppleSystemRun $d $pos "obber fortsatt med å forstå helt, men er på god vei"
$pos = Add-AppleSystemRun $d $pos "."

Write-Host "Done"
